$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data set (x, y, z) for rows 2-21: four brand new samples prepended,
# the previous last sample (old row 22) dropped.
$data = @(
    @(0.0007635815418325, 0.01328631862998, 0.0305432621389627),
    @(0.0154243474826216, -0.0059559359215199, 0.0029016099870204),
    @(0.0140499006956815, 0.0103847095742821, 0.0612392425537109),
    @(0.0195476878434419, 0.0311541277915239, 0.0459676086902618),
    @(0.0041233403608202, -0.0007635815418325, -0.066737025976181),
    @(-0.0314595587551593, 0.2884811162948608, 0.06536258012056349),
    @(-0.1331686228513717, 0.493731826543808, -0.1162171140313148),
    @(-0.09651670604944219, 0.7289149761199951, -0.0355829000473022),
    @(-0.2052507251501083, 0.5688682794570923, 0.1585195362567901),
    @(-0.18539759516716, -0.1214094683527946, -0.6217080950737),
    @(0.8633053302764893, -0.9155342578887939, 0.0097738439217209),
    @(-2.978273391723633, 1.564731359481812, -3.266449213027954),
    @(1.346346974372864, 1.957823157310486, 1.074512004852295),
    @(0.2492330223321914, 0.9819658994674684, -0.9859365224838256),
    @(0.0514653958380222, -0.3197879493236542, 0.1852448880672454),
    @(0.1069014146924018, -0.1577559560537338, 0.1151480972766876),
    @(0.7513642311096191, -1.080926060676575, -0.1253800988197326),
    @(-0.155312493443489, 0.2770273983478546, 0.052381694316864),
    @(-0.0971275717020034, -0.6624833345413208, 0.3292563557624817),
    @(-0.0768163055181503, -0.087353728711605, 0.0009162978967650999)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}

# The old sheet had 21 data rows (through row 22); the new one only has 20
# (through row 21), so drop the now-unused trailing row entirely.
$ws.Rows.Item(22).Delete()
